$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sat right after the "Instance"
#    heading run.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. The sentence "is a single and unique unit of a class. " used to be split
#    across two runs (one ending in "class" and a tiny one holding ". ").
#    Re-running Find/Replace over the full phrase merges it back into a
#    single run (keeping the formatting of the first of the two runs), which
#    is exactly what the target document looks like.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "is a single and unique unit of a class. ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "is a single and unique unit of a class. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Append a new section at the end of the document:
#      - a blank paragraph
#      - a bold "VSTS/TFS:" heading (24 half-point -> 12pt, lastRenderedPageBreak)
#      - a body paragraph whose text is split in two runs by a fresh
#        "_GoBack" bookmark (mirrors how Word leaves its own edit-resume
#        marker behind after the last place text was typed).
#    InsertXML lets us describe the exact OOXML for the new content instead
#    of relying on inherited run/paragraph formatting from whatever came
#    before it.
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>VSTS/TFS:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>VSTS and Team Foundation Server (TFS) both provide an integrated, collaborative environment that supports Git, continuous integration, and Agile tools for planning and tracking work. VSTS is the cloud offering that provides a scala</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>ble, reliable, and globally available hosted service.</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$endRange.InsertXML($newContentXml) | Out-Null
